# Update Bulgaria data spreadsheet - time_variants sheet data edits
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time_variants")

# Q2: 80 -> 50
$ws.Range("Q2").Value = 50

# Row 3: move the "0" value from H3 to G3, add a new "60" in I3, bump K3 to 95
$ws.Range("H3").Clear()
$ws.Range("G3").Value = 0
$ws.Range("I3").Value = 60
$ws.Range("K3").Value = 95

# Row 4: same pattern as row 3
$ws.Range("H4").Clear()
$ws.Range("G4").Value = 0
$ws.Range("I4").Value = 60
$ws.Range("K4").Value = 95

# S9: blank -> 0
$ws.Range("S9").Value = 0
